$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the (invisible/no-op) border formatting that was applied to A3:B6,
# restoring those cells to the default style - matches switching the
# "applyBorder with no border" style back off while keeping the real
# bottom-border style on A7:B7.
$ws.Range("A3:B6").Borders.LineStyle = -4142

# "Years" dimension (row 6): switch back to the realistic number of
# annual vehicle mileage values (3 instead of 2).
$ws.Range("B6").Value = 3

# Move the active selection to B9 (matches the saved cursor position).
$ws.Range("B9").Select() | Out-Null
